# Weekly refresh of Fruta / Hortaliza data: update D, J, K, L, M, N, O, P, Q
# columns for rows 2-10 with the new week's values (rows effectively
# reshuffled/updated in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44161; J = 50;  K = 2800; L = 3000; M = 2900; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1450; Q = 2 }
    3  = @{ D = 44225; J = 80;  K = 3400; L = 3700; M = 3550; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1775; Q = 2 }
    4  = @{ D = 44160; J = 43;  K = 3500; L = 4000; M = 3709; N = '$/paquete 36 unidades'; O = 'Región Metropolitana';  P = 103;  Q = 36 }
    5  = @{ D = 44209; J = 150; K = 3500; L = 4000; M = 3767; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1884; Q = 2 }
    6  = @{ D = 44210; J = 105; K = 3500; L = 4000; M = 3714; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1857; Q = 2 }
    7  = @{ D = 44208; J = 85;  K = 3700; L = 4000; M = 3824; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1912; Q = 2 }
    8  = @{ D = 44166; J = 70;  K = 3500; L = 4000; M = 3679; N = '$/paquete 36 unidades'; O = 'Región Metropolitana';  P = 102;  Q = 36 }
    9  = @{ D = 44215; J = 140; K = 3500; L = 4000; M = 3768; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1884; Q = 2 }
    10 = @{ D = 44223; J = 80;  K = 3500; L = 3800; M = 3688; N = '$/paquete 2 kilos';    O = 'Provincia de Diguillín'; P = 1844; Q = 2 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
}

Write-Host "Updated rows 2-10 with refreshed weekly data."
